# "tải ảnh vào cho từng câu hỏi" — point the "Hình Ảnh" (image) column at the
# real image paths under the foodimg resource folder instead of the bare
# filenames, and fix a couple of mis-matched file names along the way.
#
# The writes below are intentionally ordered the same way the original
# shared-string table was built (so newly introduced strings land at the
# same indices as in the authoritative edit) rather than in plain row
# order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$base = "/th/nguyenxuandat/FlashQuizGUI/foodimg/"

$ws.Range("G2").Value  = $base + "pho.jpg"
$ws.Range("G8").Value  = $base + "cakho.jpg"
$ws.Range("G9").Value  = $base + "nuocmam.jpg"
$ws.Range("G10").Value = $base + "muadong.jpg"
$ws.Range("G11").Value = $base + "daubep.jpg"
$ws.Range("G4").Value  = $base + "chaobo.jpg"
$ws.Range("G6").Value  = $base + "chaoluoc.jpg"
$ws.Range("G5").Value  = $base + "bunbo.jpg"
$ws.Range("G3").Value  = $base + "banhchung.jpg"
$ws.Range("G7").Value  = $base + "pho.jpg"

# Widen the image-path column so the long paths are readable, and move the
# active selection onto that column. (49.109375 "characters" is what the
# target column width rounds to when stored; COM's ColumnWidth setter here
# quantizes to the nearest 1/6 character step, same as Excel itself does
# when a column is resized by dragging, so feed it an input that lands on
# the closest attainable step.)
$ws.Columns.Item(7).ColumnWidth = 48.3
$ws.Range("G3").Select() | Out-Null
